# Update Lama1-Itga7 LR-pairs data with new TPM-based values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 1).Value = "FAPs"
$ws.Cells.Item(2, 2).Value = "Lama1"
$ws.Cells.Item(2, 3).Value = "Itga7"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.5587383333333333
$ws.Cells.Item(2, 8).Value = 1.676215
$ws.Cells.Item(2, 9).Value = 0.8486764927018626
$ws.Cells.Item(2, 10).Value = 0.8937587278261895
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 2.035934
$ws.Cells.Item(2, 14).Value = 6.107802
$ws.Cells.Item(2, 15).Value = 0.03126880699186227
$ws.Cells.Item(2, 16).Value = 0.04430738339814538
$ws.Cells.Item(2, 17).Value = 1.137554369936667
$ws.Cells.Item(2, 18).Value = 10.23798932943
$ws.Cells.Item(2, 19).Value = 0.02653710144882515
$ws.Cells.Item(2, 20).Value = 0.03960011061923364

# Row 3
$ws.Cells.Item(3, 1).Value = "FAPs"
$ws.Cells.Item(3, 2).Value = "Lama1"
$ws.Cells.Item(3, 3).Value = "Itga7"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.5587383333333333
$ws.Cells.Item(3, 8).Value = 1.676215
$ws.Cells.Item(3, 9).Value = 0.8486764927018626
$ws.Cells.Item(3, 10).Value = 0.8937587278261895
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 3.794922333333334
$ws.Cells.Item(3, 14).Value = 11.384767
$ws.Cells.Item(3, 15).Value = 0.0582841555718936
$ws.Cells.Item(3, 16).Value = 0.08258768643246023
$ws.Cells.Item(3, 17).Value = 2.120368579656111
$ws.Cells.Item(3, 18).Value = 19.083317216905
$ws.Cells.Item(3, 19).Value = 0.04946439273084438
$ws.Cells.Item(3, 20).Value = 0.0738134655599839

# Row 4
$ws.Cells.Item(4, 1).Value = "FAPs"
$ws.Cells.Item(4, 2).Value = "Lama1"
$ws.Cells.Item(4, 3).Value = "Itga7"
$ws.Cells.Item(4, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.5587383333333333
$ws.Cells.Item(4, 8).Value = 1.676215
$ws.Cells.Item(4, 9).Value = 0.8486764927018626
$ws.Cells.Item(4, 10).Value = 0.8937587278261895
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 0.451002
$ws.Cells.Item(4, 14).Value = 1.353006
$ws.Cells.Item(4, 15).Value = 0.006926695310822388
$ws.Cells.Item(4, 16).Value = 0.009815012926416261
$ws.Cells.Item(4, 17).Value = 0.25199210581
$ws.Cells.Item(4, 18).Value = 2.26792895229
$ws.Cells.Item(4, 19).Value = 0.005878523482403183
$ws.Cells.Item(4, 20).Value = 0.008772253466711403

# Row 5
$ws.Cells.Item(5, 1).Value = "FAPs"
$ws.Cells.Item(5, 2).Value = "Lama1"
$ws.Cells.Item(5, 3).Value = "Itga7"
$ws.Cells.Item(5, 4).Value = "MuSCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.5587383333333333
$ws.Cells.Item(5, 8).Value = 1.676215
$ws.Cells.Item(5, 9).Value = 0.8486764927018626
$ws.Cells.Item(5, 10).Value = 0.8937587278261895
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 57.4814495
$ws.Cells.Item(5, 14).Value = 114.962899
$ws.Cells.Item(5, 15).Value = 0.8828264325012393
$ws.Cells.Item(5, 16).Value = 0.8339669888701803
$ws.Cells.Item(5, 17).Value = 32.11708929121416
$ws.Cells.Item(5, 18).Value = 192.702535747285
$ws.Cells.Item(5, 19).Value = 0.7492340403996494
$ws.Cells.Item(5, 20).Value = 0.7453652750216503

# Row 6
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Lama1"
$ws.Cells.Item(6, 3).Value = "Itga7"
$ws.Cells.Item(6, 4).Value = "Resolving-Mac"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 0.5587383333333333
$ws.Cells.Item(6, 8).Value = 1.676215
$ws.Cells.Item(6, 9).Value = 0.8486764927018626
$ws.Cells.Item(6, 10).Value = 0.8937587278261895
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 1.347395
$ws.Cells.Item(6, 14).Value = 4.042185
$ws.Cells.Item(6, 15).Value = 0.02069390962418245
$ws.Cells.Item(6, 16).Value = 0.02932292837279799
$ws.Cells.Item(6, 17).Value = 0.7528412366416666
$ws.Cells.Item(6, 18).Value = 6.775571129775
$ws.Cells.Item(6, 19).Value = 0.01756243464014048
$ws.Cells.Item(6, 20).Value = 0.02620762315861041

# Row 7
$ws.Cells.Item(7, 1).Value = "MuSCs"
$ws.Cells.Item(7, 2).Value = "Lama1"
$ws.Cells.Item(7, 3).Value = "Itga7"
$ws.Cells.Item(7, 4).Value = "Inflammatory-Mac"
$ws.Cells.Item(7, 5).Value = 2
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 0.09962599999999999
$ws.Cells.Item(7, 8).Value = 0.199252
$ws.Cells.Item(7, 9).Value = 0.1513235072981373
$ws.Cells.Item(7, 10).Value = 0.1062412721738106
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 2.035934
$ws.Cells.Item(7, 14).Value = 6.107802
$ws.Cells.Item(7, 15).Value = 0.03126880699186227
$ws.Cells.Item(7, 16).Value = 0.04430738339814538
$ws.Cells.Item(7, 17).Value = 0.2028319606839999
$ws.Cells.Item(7, 18).Value = 1.216991764104
$ws.Cells.Item(7, 19).Value = 0.004731705543037118
$ws.Cells.Item(7, 20).Value = 0.00470727277891174

# Row 8
$ws.Cells.Item(8, 1).Value = "MuSCs"
$ws.Cells.Item(8, 2).Value = "Lama1"
$ws.Cells.Item(8, 3).Value = "Itga7"
$ws.Cells.Item(8, 4).Value = "FAPs"
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 0.09962599999999999
$ws.Cells.Item(8, 8).Value = 0.199252
$ws.Cells.Item(8, 9).Value = 0.1513235072981373
$ws.Cells.Item(8, 10).Value = 0.1062412721738106
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 3.794922333333334
$ws.Cells.Item(8, 14).Value = 11.384767
$ws.Cells.Item(8, 15).Value = 0.0582841555718936
$ws.Cells.Item(8, 16).Value = 0.08258768643246023
$ws.Cells.Item(8, 17).Value = 0.3780729323806666
$ws.Cells.Item(8, 18).Value = 2.268437594284
$ws.Cells.Item(8, 19).Value = 0.008819762841049213
$ws.Cells.Item(8, 20).Value = 0.00877422087247633

# Row 9
$ws.Cells.Item(9, 1).Value = "MuSCs"
$ws.Cells.Item(9, 2).Value = "Lama1"
$ws.Cells.Item(9, 3).Value = "Itga7"
$ws.Cells.Item(9, 4).Value = "MuSCs"
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 0.09962599999999999
$ws.Cells.Item(9, 8).Value = 0.199252
$ws.Cells.Item(9, 9).Value = 0.1513235072981373
$ws.Cells.Item(9, 10).Value = 0.1062412721738106
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.451002
$ws.Cells.Item(9, 14).Value = 1.353006
$ws.Cells.Item(9, 15).Value = 0.006926695310822388
$ws.Cells.Item(9, 16).Value = 0.009815012926416261
$ws.Cells.Item(9, 17).Value = 0.04493152525199999
$ws.Cells.Item(9, 18).Value = 0.269589151512
$ws.Cells.Item(9, 19).Value = 0.001048171828419205
$ws.Cells.Item(9, 20).Value = 0.001042759459704859

# Row 10
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Lama1"
$ws.Cells.Item(10, 3).Value = "Itga7"
$ws.Cells.Item(10, 4).Value = "MuSCs"
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.09962599999999999
$ws.Cells.Item(10, 8).Value = 0.199252
$ws.Cells.Item(10, 9).Value = 0.1513235072981373
$ws.Cells.Item(10, 10).Value = 0.1062412721738106
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 57.4814495
$ws.Cells.Item(10, 14).Value = 114.962899
$ws.Cells.Item(10, 15).Value = 0.8828264325012393
$ws.Cells.Item(10, 16).Value = 0.8339669888701803
$ws.Cells.Item(10, 17).Value = 5.726646887886999
$ws.Cells.Item(10, 18).Value = 22.906587551548
$ws.Cells.Item(10, 19).Value = 0.1335923921015898
$ws.Cells.Item(10, 20).Value = 0.08860171384853008

# Row 11
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Lama1"
$ws.Cells.Item(11, 3).Value = "Itga7"
$ws.Cells.Item(11, 4).Value = "Resolving-Mac"
$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.09962599999999999
$ws.Cells.Item(11, 8).Value = 0.199252
$ws.Cells.Item(11, 9).Value = 0.1513235072981373
$ws.Cells.Item(11, 10).Value = 0.1062412721738106
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 1.347395
$ws.Cells.Item(11, 14).Value = 4.042185
$ws.Cells.Item(11, 15).Value = 0.02069390962418245
$ws.Cells.Item(11, 16).Value = 0.02932292837279799
$ws.Cells.Item(11, 17).Value = 0.13423557427
$ws.Cells.Item(11, 18).Value = 0.80541344562
$ws.Cells.Item(11, 19).Value = 0.003131474984041967
$ws.Cells.Item(11, 20).Value = 0.003115305214187584

